$d = $word.ActiveDocument

$d.Content.Find.Execute("437÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "530÷9=", 2) | Out-Null
$d.Content.Find.Execute("794÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "863÷6=", 2) | Out-Null
$d.Content.Find.Execute("384÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "414÷6=", 2) | Out-Null
$d.Content.Find.Execute("669÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "937÷2=", 2) | Out-Null
$d.Content.Find.Execute("827÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "589÷5=", 2) | Out-Null
$d.Content.Find.Execute("423÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "995÷7=", 2) | Out-Null
$d.Content.Find.Execute("854÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "859÷8=", 2) | Out-Null
$d.Content.Find.Execute("476÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "160÷8=", 2) | Out-Null
$d.Content.Find.Execute("580÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "598÷7=", 2) | Out-Null
$d.Content.Find.Execute("789÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "691÷4=", 2) | Out-Null
$d.Content.Find.Execute("462÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "216÷4=", 2) | Out-Null
$d.Content.Find.Execute("518÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "179÷5=", 2) | Out-Null
$d.Content.Find.Execute("451÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "219÷4=", 2) | Out-Null
$d.Content.Find.Execute("160÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "277÷4=", 2) | Out-Null
$d.Content.Find.Execute("197÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "628÷6=", 2) | Out-Null
$d.Content.Find.Execute("890÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "124÷5=", 2) | Out-Null
$d.Content.Find.Execute("188÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "612÷9=", 2) | Out-Null
$d.Content.Find.Execute("434÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "257÷4=", 2) | Out-Null
$d.Content.Find.Execute("450÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "133÷7=", 2) | Out-Null
$d.Content.Find.Execute("815÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "618÷4=", 2) | Out-Null
$d.Content.Find.Execute("841÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "328÷3=", 2) | Out-Null
$d.Content.Find.Execute("499÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "556÷7=", 2) | Out-Null
$d.Content.Find.Execute("395÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "906÷2=", 2) | Out-Null
$d.Content.Find.Execute("553÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "839÷8=", 2) | Out-Null
$d.Content.Find.Execute("680÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "213÷7=", 2) | Out-Null
